$d = $word.ActiveDocument

# 1. "Native American" + " Address" (two runs) -> "American Indian Address Part " + "2"
#    Locate the two runs via Find (process the later one first so the earlier
#    run's character offsets stay valid).
$rngAddress = $d.Content
$rngAddress.Find.Execute(" Address", $true, $false, $false, $false, $false,
                          $true, 1, $false, "", 0)
$addressStart = $rngAddress.Start
$addressEnd = $rngAddress.End

$rngNative = $d.Content
$rngNative.Find.Execute("Native American", $true, $false, $false, $false, $false,
                         $true, 1, $false, "", 0)
$nativeStart = $rngNative.Start
$nativeEnd = $rngNative.End

$runAddress = $d.Range($addressStart, $addressEnd)
$runAddress.Text = "2"
$runAddress.Font.Name = "Times New Roman"

$runNative = $d.Range($nativeStart, $nativeEnd)
$runNative.Text = "American Indian Address Part "
$runNative.Font.Name = "Times New Roman"

# 2. "Atrial Fibrillation" -> "American Indian Address Part 1"
$d.Content.Find.Execute("Atrial Fibrillation", $true, $false, $false, $false, $false,
                         $true, 1, $false, "American Indian Address Part 1", 2)

# 3. "M 12/13" -> "TBD"
$d.Content.Find.Execute("M 12/13", $true, $false, $false, $false, $false,
                         $true, 1, $false, "TBD", 2)

Write-Output "done"
